$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-12-14"

# Update the December label cell to reflect the new "through" date
$ws.Range("A13").Value = "December (through 12-14)"

# Update December row (row 13) values for years 2015-2022 (columns B-I)
$ws.Range("B13").Value = 15
$ws.Range("C13").Value = 44
$ws.Range("D13").Value = 46
$ws.Range("E13").Value = 31
$ws.Range("F13").Value = 23
$ws.Range("G13").Value = 69
$ws.Range("H13").Value = 103
$ws.Range("I13").Value = 60

# Update Total row (row 14) values for years 2015-2022 (columns B-I)
$ws.Range("B14").Value = 306
$ws.Range("C14").Value = 607
$ws.Range("D14").Value = 867
$ws.Range("E14").Value = 713
$ws.Range("F14").Value = 557
$ws.Range("G14").Value = 1333
$ws.Range("H14").Value = 1746
$ws.Range("I14").Value = 1576
